# Tugas 1 (Implementasi Import pada Semua Tabel)
# Fix the typo in the template's header: "kateori_id" -> "kategori_id".
# (All other data in the sheet is unchanged; the shared-string table
# reorders itself automatically when the edited string is re-saved.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "kategori_id"

# Restore the cursor/selection position recorded in the saved file.
[void]$ws.Range("G9").Select()
